# Fix - Purchase service tests
# Remove the two extra product rows (Widget C, Widget E) that were
# added by mistake, and correct the Quantity for row 3 (Widget D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct Quantity value for Widget D (row 3, column E) 26 -> 25
$ws.Range("E3").Value = 25

# Delete rows 6 and 5 (delete higher row first to keep indices stable)
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
